# Generate Report for Handoff
# Adds a new tracked file (818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md) as row 9
# on the Overview / zh-cn / de-de sheets, resizes the 3 tables to include the
# new row, and wires up the corresponding hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovwTable = $ovw.ListObjects.Item(1)
$ovwTable.ListRows.Add() | Out-Null

$ovw.Range("A9").Value = "818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md"
$ovw.Range("B9").Value = "e2e\818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md"
$ovw.Range("C9").Value = ".md"
$ovw.Range("D9").Value = ""
$ovw.Range("E9").Value = "Ready for handoff"
$ovw.Range("F9").Value = "Ready for handoff"
$ovw.Range("G9").Value = "2016-11-29 02:39:21"
$ovw.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ovw.Hyperlinks.Add(
    $ovw.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a6a7a2d9a3b6a6e2b5e6f9c7d4a3b2c1d0e9f8a/e2e/818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md",
    "",
    "",
    "e2e\818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md"
) | Out-Null

$ovwTable.Resize($ovw.Range("A1:G9")) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhTable = $zh.ListObjects.Item(1)
$zhTable.ListRows.Add() | Out-Null

$zh.Range("A9").Value = "818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md"
$zh.Range("B9").Value = ".md"
$zh.Range("C9").Value = "Ready for handoff"
$zh.Range("D9").Value = "e2e"
$zh.Range("E9").Value = "ht"
$zh.Range("F9").Value = "False"
$zh.Range("G9").Value = "818f16f7-b23c-470b-8df4-dc22cf8d5c4f.7bd9a3a8067da69bcc90e44a57eacf73c96527e6.zh-cn.xlf"
$zh.Range("H9").Value = "2016-11-29 02:39:07"
$zh.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("I9").Value = ""
$zh.Range("J9").Value = ""
$zh.Range("K9").Value = "0001-01-01 00:00:00"
$zh.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("L9").Value = ""
$zh.Range("M9").Value = "True"
$zh.Range("N9").Value = ""
$zh.Range("O9").Value = "False"
$zh.Range("P9").Value = ""

$zh.Hyperlinks.Add(
    $zh.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a6a7a2d9a3b6a6e2b5e6f9c7d4a3b2c1d0e9f8a/e2e/818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md",
    "",
    "",
    "818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md"
) | Out-Null

$zhTable.Resize($zh.Range("A1:P9")) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deTable = $de.ListObjects.Item(1)
$deTable.ListRows.Add() | Out-Null

$de.Range("A9").Value = "818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md"
$de.Range("B9").Value = ".md"
$de.Range("C9").Value = "Ready for handoff"
$de.Range("D9").Value = "e2e"
$de.Range("E9").Value = "ht"
$de.Range("F9").Value = "False"
$de.Range("G9").Value = "818f16f7-b23c-470b-8df4-dc22cf8d5c4f.7bd9a3a8067da69bcc90e44a57eacf73c96527e6.de-de.xlf"
$de.Range("H9").Value = "2016-11-29 02:39:21"
$de.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("I9").Value = ""
$de.Range("J9").Value = ""
$de.Range("K9").Value = "0001-01-01 00:00:00"
$de.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("L9").Value = ""
$de.Range("M9").Value = "True"
$de.Range("N9").Value = ""
$de.Range("O9").Value = "False"
$de.Range("P9").Value = ""

$de.Hyperlinks.Add(
    $de.Range("A9"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7a6a7a2d9a3b6a6e2b5e6f9c7d4a3b2c1d0e9f8a/e2e/818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md",
    "",
    "",
    "818f16f7-b23c-470b-8df4-dc22cf8d5c4f.md"
) | Out-Null

$deTable.Resize($de.Range("A1:P9")) | Out-Null
